# Auto-generated edit script: updates cryptocurrency price/volume data
# and swaps the PEPE/Aptos rows (44 <-> 45) per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.845.25"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "'1.916.68"
$ws.Range("E3").Value = "  +0.72%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'324.04"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").Value = "'0.4560"
$ws.Range("E7").Value = "  -0.76%  "

$ws.Range("D8").Value = "'0.3803"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.07749"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("D10").Value = "'0.9760"
$ws.Range("E10").Value = "  -0.40%  "

$ws.Range("D11").Value = "'22.20"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").Value = "'1.935.47"
$ws.Range("E12").Value = "  +3.74%  "

$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").Value = "'0.06986"
$ws.Range("E15").Value = "  -0.76%  "

$ws.Range("D17").Value = "'84.36"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("D18").Value = "'0.000009490"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("D19").Value = "'16.63"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D21").Value = "'28.864.37"
$ws.Range("E21").Value = "  -0.24%  "

$ws.Range("D22").Value = "'5.341"
$ws.Range("E22").Value = "  +0.52%  "

$ws.Range("D23").Value = "'11.12"
$ws.Range("E23").Value = "  +2.22%  "

$ws.Range("D24").Value = "'2.154.48"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("D26").Value = "'157.96"
$ws.Range("E26").Value = "  +0.80%  "

$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").Value = "'5.615"
$ws.Range("E28").Value = "  +0.51%  "

$ws.Range("D29").Value = "'117.71"
$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("D30").Value = "'1.854"
$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("D31").Value = "'0.09280"
$ws.Range("E31").Value = "  +0.21%  "

$ws.Range("D32").Value = "'0.8698"
$ws.Range("E32").Value = "  +1.24%  "

$ws.Range("D33").Value = "'5.104"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("D34").Value = "'1.245"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("E35").Value = "  +1.01%  "

$ws.Range("D36").Value = "'0.05700"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("D37").Value = "'1.149"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").Value = "'0.02037"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("E40").Value = "  +11.06%  "

$ws.Range("D41").Value = "'7.503"
$ws.Range("E41").Value = "  +0.63%  "

$ws.Range("D42").Value = "'0.5501"

$ws.Range("D43").Value = "'0.1757"
$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D46").Value = "'2.168"
$ws.Range("E46").Value = "  +3.65%  "

$ws.Range("D47").Value = "'0.5160"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").Value = "'0.06928"
$ws.Range("E48").Value = "  +1.80%  "

$ws.Range("D49").Value = "'11.13"
$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").Value = "'110.56"
$ws.Range("E50").Value = "  -0.65%  "

$ws.Range("D51").Value = "'1.762"
$ws.Range("E51").Value = "  -0.83%  "

# Row 44/45: PEPE and Aptos swap positions with updated price/volume data
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'9.337"
$ws.Range("E44").Value = "  +0.94%  "

$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "'0.000002870"
$ws.Range("E45").Value = "  +16.22%  "
